$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''63.016.46'
$ws.Range('E2').Value = '  +2.89%  '
$ws.Range('D3').Value = '''2.467.69'
$ws.Range('E3').Value = '  +4.71%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '''566.23'
$ws.Range('E5').Value = '  +1.52%  '
$ws.Range('D6').Value = '''142.60'
$ws.Range('E6').Value = '  +7.23%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E8').Value = '  +0.49%  '
$ws.Range('D9').Value = '''2.467.09'
$ws.Range('E9').Value = '  +4.85%  '
$ws.Range('E10').Value = '  +1.97%  '
$ws.Range('E11').Value = '  +0.84%  '
$ws.Range('E12').Value = '  +1.22%  '
$ws.Range('E13').Value = '  +2.71%  '
$ws.Range('D14').Value = '''26.64'
$ws.Range('E14').Value = '  +9.13%  '
$ws.Range('D15').Value = '''2.910.15'
$ws.Range('E15').Value = '  +4.72%  '
$ws.Range('D16').Value = '''62.887.20'
$ws.Range('E16').Value = '  +2.81%  '
$ws.Range('E17').Value = '  +3.71%  '
$ws.Range('D18').Value = '''2.468.72'
$ws.Range('E18').Value = '  +4.40%  '
$ws.Range('E19').Value = '  +3.87%  '
$ws.Range('D20').Value = '''340.23'
$ws.Range('E20').Value = '  +6.86%  '
$ws.Range('E21').Value = '  +2.60%  '
$ws.Range('D22').Value = '''6.82'
$ws.Range('E22').Value = '  +1.60%  '
$ws.Range('E23').Value = '  +0.15%  '
$ws.Range('D24').Value = '''65.57'
$ws.Range('E24').Value = '  +1.44%  '
$ws.Range('D25').Value = '''0.172'
$ws.Range('E25').Value = '  -0.60%  '
$ws.Range('D26').Value = '''0.999'
$ws.Range('E26').Value = '  +0.08%  '
$ws.Range('E27').Value = '  +4.38%  '
$ws.Range('D28').Value = '''8.09'
$ws.Range('E28').Value = '  -0.09%  '
$ws.Range('D29').Value = '''1.39'
$ws.Range('E29').Value = '  +7.96%  '
$ws.Range('D30').Value = '''6.82'
$ws.Range('E30').Value = '  +10.34%  '
$ws.Range('E31').Value = '  +5.01%  '
$ws.Range('D32').Value = '''0.0₃0798'
$ws.Range('E32').Value = '  +6.96%  '
$ws.Range('D33').Value = '''175.04'
$ws.Range('E33').Value = '  +2.33%  '
$ws.Range('E34').Value = '  +8.78%  '
$ws.Range('E35').Value = '  +2.66%  '
$ws.Range('D36').Value = '''18.80'
$ws.Range('E36').Value = '  +3.14%  '
$ws.Range('D37').Value = '''374.06'
$ws.Range('E37').Value = '  +12.24%  '
$ws.Range('E38').Value = '  +0.02%  '
$ws.Range('E39').Value = '  +3.88%  '
$ws.Range('D40').Value = '''0.999'
$ws.Range('E40').Value = '  -0.16%  '
$ws.Range('E41').Value = '  +9.08%  '
$ws.Range('D42').Value = '''40.32'
$ws.Range('D43').Value = '''150.21'
$ws.Range('E43').Value = '  +7.12%  '
$ws.Range('D44').Value = '''3.69'
$ws.Range('E44').Value = '  +3.45%  '
$ws.Range('E45').Value = '  +5.37%  '
$ws.Range('E46').Value = '  +4.49%  '
$ws.Range('E47').Value = '  +0.31%  '
$ws.Range('E48').Value = '  +2.31%  '
$ws.Range('E49').Value = '  +3.75%  '
$ws.Range('E50').Value = '  +2.53%  '
$ws.Range('D51').Value = '''17.89'
$ws.Range('E51').Value = '  +2.54%  '
